# Applies updated crypto price/volume data (and one row insertion for
# WrappedliquidstakedEther2.0) to match the refreshed GitHub Actions feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.485.02'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '1.845.84'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '''261.12'
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = '''0.5243'
$ws.Range("E7").Value = '  +1.47%  '
$ws.Range("D8").Value = '''0.3213'
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").Value = '''0.06782'
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("D10").Value = '''18.81'
$ws.Range("E10").Value = '  +0.33%  '
$ws.Range("D11").Value = '''0.7792'
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("D12").Value = '''0.07751'
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("D13").Value = '1.841.61'
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("D14").Value = '''88.00'
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").Value = '''5.006'
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").Value = '''13.90'
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").Value = '''0.000007927'
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("D20").Value = '26.513.29'
$ws.Range("E20").Value = '  +0.33%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.081.35'
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '''4.613'
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '''9.418'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '''5.970'
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '''141.73'
$ws.Range("E25").Value = '  -1.97%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '''2.150'
$ws.Range("E26").Value = '  -6.13%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '''1.675'
$ws.Range("E27").Value = '  +2.52%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '''16.95'
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '''111.71'
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '''4.166'
$ws.Range("E30").Value = '  -0.80%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '''0.08708'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''4.062'
$ws.Range("E32").Value = '  -1.57%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.04861'
$ws.Range("E33").Value = '  +0.65%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''0.7192'
$ws.Range("E34").Value = '  +4.13%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.873'
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''1.126'
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = '''3.100'
$ws.Range("E37").Value = '  +0.42%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '''2.252'
$ws.Range("E38").Value = '  +2.28%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.01771'
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '''0.4817'
$ws.Range("E40").Value = '  -2.05%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '''0.8946'
$ws.Range("E41").Value = '  +0.56%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '''110.15'
$ws.Range("E42").Value = '  -1.68%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.925'
$ws.Range("E43").Value = '  -3.42%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '''1.001'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '''7.648'
$ws.Range("E45").Value = '  -1.38%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.4149'
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.05872'
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''8.929'
$ws.Range("E48").Value = '  -1.78%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''34.97'
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '''0.1229'
$ws.Range("E50").Value = '  -1.91%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = '''0.8888'
$ws.Range("E51").Value = '  +0.80%  '
